# v1.1 refactor: Parameters sheet now drives a data.table-based pipeline,
# so the old MCMC-model tuning knobs (model_fc, norm_*, model_*, nworker)
# collapse down to just two new knobs: max_spectra_per_peptide, nbatch.

$wb = $excel.ActiveWorkbook

$wsParameters = $wb.Worksheets.Item("Parameters")
$wsDesign     = $wb.Worksheets.Item("Design")
$wsFractions  = $wb.Worksheets.Item("Fractions")

# Drop the old rows 4-11 (norm_nburnin .. nworker) entirely.
$wsParameters.Range("A4:B11").EntireRow.Delete()

# Relabel / re-value the remaining two parameter rows.
$wsParameters.Range("A2").Value = "max_spectra_per_peptide"
$wsParameters.Range("B2").Value = 5
$wsParameters.Range("A3").Value = "nbatch"
$wsParameters.Range("B3").Value = 100

# Restore the selections on the other two sheets ...
[void]$wsDesign.Range("C12").Select()
[void]$wsFractions.Range("F12").Select()

# ... then make Parameters the active tab with its own selection, last,
# so it ends up the active sheet in the saved workbook.
[void]$wsParameters.Activate()
[void]$wsParameters.Range("C19").Select()
